$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (coin name) and C (link) swaps for rows 33/34
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

# Column D (price) updates - force text so values like "23.461.47" or "1.002" are not coerced to numbers/dates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.461.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3813"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3616"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08273"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.233"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.469"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.351"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.638.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06964"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.601"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.474.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.531"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.076"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.272"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.819.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.091"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.156"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.565"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02772"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2516"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08772"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.988"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07028"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7054"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.348"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6540"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.296"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.964"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07982"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.45"
$ws.Range("D50").Style = "Normal"

# Column E (volume/percent change) updates
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  +4.15%  "
$ws.Range("E26").Value = "  -4.27%  "
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").Value = "  +15.41%  "
$ws.Range("E33").Value = "  -6.56%  "
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("E35").Value = "  +5.96%  "
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("E51").Value = "  -0.85%  "
